$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.425.24"
$ws.Range("E2").Value = "  -0.29%  "

# Row 3
$ws.Range("D3").Value = "3.324.14"
$ws.Range("E3").Value = "  +0.63%  "

# Row 4
$ws.Range("E4").Value = "  -0.40%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.55"
$ws.Range("E5").Value = "  -0.68%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.90"
$ws.Range("E6").Value = "  -3.79%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.34%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.588"
$ws.Range("E8").Value = "  +0.13%  "

# Row 9
$ws.Range("D9").Value = "3.322.59"
$ws.Range("E9").Value = "  +1.05%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.178"
$ws.Range("E10").Value = "  -0.11%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.576"
$ws.Range("E11").Value = "  -0.08%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.48"
$ws.Range("E12").Value = "  -1.85%  "

# Row 13
$ws.Range("E13").Value = "  -2.07%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "660.17"
$ws.Range("E14").Value = "  +3.78%  "

# Row 15
$ws.Range("D15").Value = "3.861.84"
$ws.Range("E15").Value = "  +0.33%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.41"
$ws.Range("E16").Value = "  +0.08%  "

# Row 17
$ws.Range("D17").Value = "67.555.19"
$ws.Range("E17").Value = "  -0.55%  "

# Row 18
$ws.Range("E18").Value = "  -0.35%  "

# Row 19
$ws.Range("D19").Value = "3.326.60"
$ws.Range("E19").Value = "  -0.60%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.40"
$ws.Range("E20").Value = "  -1.13%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.95"
$ws.Range("E21").Value = "  +0.51%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.889"
$ws.Range("E22").Value = "  -0.68%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.45"
$ws.Range("E23").Value = "  +8.38%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.06"
$ws.Range("E24").Value = "  -3.55%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "99.61"
$ws.Range("E25").Value = "  +2.06%  "

# Row 26
$ws.Range("E26").Value = "  -3.72%  "

# Row 27
$ws.Range("E27").Value = "  -4.67%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.27"
$ws.Range("E28").Value = "  -3.16%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.71"
$ws.Range("E29").Value = "  +3.22%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.46"
$ws.Range("E30").Value = "  +11.61%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.43"
$ws.Range("E31").Value = "  -1.17%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "578.08"
$ws.Range("E32").Value = "  -3.15%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "10.97"
$ws.Range("E33").Value = "  +0.54%  "

# Row 34
$ws.Range("E34").Value = "  +0.15%  "

# Row 35
$ws.Range("E35").Value = "  -1.88%  "

# Row 36
$ws.Range("D36").Value = "3.695.56"
$ws.Range("E36").Value = "  -6.43%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.61"
$ws.Range("E37").Value = "  +1.61%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.35"
$ws.Range("E38").Value = "  -7.43%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.54"
$ws.Range("E39").Value = "  +4.92%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.131"
$ws.Range("E40").Value = "  +1.49%  "

# Row 41
$ws.Range("E41").Value = "  -3.35%  "

# Row 42
$ws.Range("E42").Value = "  -4.17%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.333"
$ws.Range("E43").Value = "  -0.91%  "

# Row 44
$ws.Range("D44").Value = "0.0₃0666"
$ws.Range("E44").Value = "  -3.40%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.28"
$ws.Range("E45").Value = "  -3.24%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0406"
$ws.Range("E46").Value = "  -2.30%  "

# Row 47
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.59"
$ws.Range("E47").Value = "  +1.91%  "

# Row 48
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.128"
$ws.Range("E48").Value = "  -0.03%  "

# Row 49
$ws.Range("E49").Value = "  -0.59%  "

# Row 50
$ws.Range("E50").Value = "  +0.99%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "127.42"
$ws.Range("E51").Value = "  -2.48%  "
